# Curated ICDC startup workbook (TC47 Yorkshire Terrier): refresh the
# cart/case/sample/file Cypher queries, add a StudyFilesTab row, and
# refresh the row/column layout to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: D1 (cartQuery label) ---------------------------------------
$ws.Range("D1").Value2 = 'cartQuery'

# --- Row 2: CasesTab ----------------------------------------------------------
# (write B2 first so its text claims the next shared-string slot, matching
#  the authored workbooks string order)
$ws.Range("B2").Value2 = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN [''Yorkshire Terrier'']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age
RETURN  
       coalesce(c.case_id, '''') AS `Case ID`,
       coalesce(s.clinical_study_designation, '''') AS `Study Code`,
       coalesce(s.clinical_study_type, '''') AS  `Study Type`,
       coalesce(demo.breed, '''') AS Breed ,
       coalesce(diag.disease_term, '''') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '''') AS `Stage of Disease`,
       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,
       coalesce(demo.sex, '''') AS Sex,
       coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
       coalesce(demo.weight, '''') AS `Weight (kg)`,
       coalesce(diag.best_response, '''') AS `Response to Treatment`,
       coalesce(co.cohort_description, '''') AS `Cohort`'

# --- Row 5: StudyFilesTab (new row) --------------------------------------------
$ws.Range("A5").Value2 = 'StudyFilesTab'
$ws.Range("B5").Value2 = '  MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
WHERE demo.breed  IN [''Yorkshire Terrier''] 
WITH DISTINCT f,  s, c, demo, diag
WITH
        f, c, demo, diag, s,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH    
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH    
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '''') AS `File Name`,
  coalesce(f.file_type, '''') AS `File Type`,
  coalesce("study", '''') AS `Association`,
  coalesce(f.file_description, '''') AS `Description`,
  coalesce(f.file_format, '''') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
  coalesce(s.clinical_study_designation,'''') AS `Study Code`'

# --- Row 3: SamplesTab ----------------------------------------------------------
$ws.Range("B3").Value2 = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE demo.breed IN [''Yorkshire Terrier''] 
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '''') AS `Sample ID`, 
        coalesce(c.case_id, '''') AS `Case ID`, 
        coalesce(demo.breed,'''') AS Breed,
        coalesce(diag.disease_term,'''') AS Diagnosis, 
        coalesce(samp.sample_site, '''') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '''') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '''') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '''') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '''') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '''') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '''') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '''') AS `Sample Preservation`'

# --- Row 4: FilesTab -------------------------------------------------------------
$ws.Range("B4").Value2 = 'MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f)-[*]->(samp:sample)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN [''Yorkshire Terrier''] 
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN 
        coalesce(f.file_name, '''') AS `File Name`,
        coalesce(f.file_type, '''') AS `File Type`,
        coalesce(labels(parent)[0], '''') AS `Association`,
        coalesce(f.file_description, '''') AS `Description`,
        coalesce(f.file_format, '''') AS `Format`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
        coalesce(samp.sample_id, '''') AS `Sample ID`,
        coalesce(c.case_id, '''') AS `Case ID`,
        coalesce(demo.breed,'''') AS Breed ,
        coalesce(diag.disease_term,'''') AS Diagnosis'

# --- StatQuery / cartQuery columns (C,D) now share the same text on every row ---
$cartQueryText = 'MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN [''Yorkshire Terrier'']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`'
$ws.Range("C2").Value2 = $cartQueryText
$ws.Range("D2").Value2 = $cartQueryText
$ws.Range("C3").Value2 = $cartQueryText
$ws.Range("D3").Value2 = $cartQueryText
$ws.Range("C4").Value2 = $cartQueryText
$ws.Range("D4").Value2 = $cartQueryText
$ws.Range("C5").Value2 = $cartQueryText
$ws.Range("D5").Value2 = $cartQueryText

# --- Row 5 output filename columns match the other rows -------------------------
$ws.Range("E5").Value2 = $ws.Range("E4").Value2
$ws.Range("F5").Value2 = $ws.Range("F4").Value2

# --- Wrap text on B:D for every data row (reuses the existing wrap-text style) --
$ws.Range("B2:D5").WrapText = $true

# --- Row heights ------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 288
$ws.Rows.Item(3).RowHeight = 145.5
$ws.Rows.Item(4).RowHeight = 102.75
$ws.Rows.Item(5).RowHeight = 102.75

# --- Column widths ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.85546875
$ws.Columns.Item(2).ColumnWidth = 92.42578125
$ws.Columns.Item(3).ColumnWidth = 75.7109375
$ws.Columns.Item(4).ColumnWidth = 75.7109375
$ws.Columns.Item(5).ColumnWidth = 70.28515625
$ws.Columns.Item(6).ColumnWidth = 28.5703125

# --- View / selection (matches the saved workbook state) --------------------------
$ws.Range("A3").Select()
$ws.Range("B9").Select()

